$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "29.619.56"
Set-TextValue "E2" "  +4.43%  "
Set-TextValue "D3" "1.603.70"
Set-TextValue "E3" "  +3.42%  "
Set-TextValue "E4" "  -0.51%  "
Set-TextValue "D5" "213.60"
Set-TextValue "E5" "  +1.84%  "
Set-TextValue "D6" "0.517"
Set-TextValue "E6" "  +7.40%  "
Set-TextValue "E7" "  -0.51%  "
Set-TextValue "D8" "26.82"
Set-TextValue "E8" "  +12.20%  "
Set-TextValue "E9" "  +3.51%  "
Set-TextValue "E10" "  +3.05%  "
Set-TextValue "D11" "0.0915"
Set-TextValue "E11" "  +2.91%  "
Set-TextValue "D12" "1.831.55"
Set-TextValue "E12" "  +3.42%  "
Set-TextValue "D13" "1.595.78"
Set-TextValue "E13" "  +2.86%  "
Set-TextValue "D14" "29.657.62"
Set-TextValue "E14" "  +4.66%  "
Set-TextValue "D15" "3.77"
Set-TextValue "E15" "  +3.95%  "
Set-TextValue "E16" "  +3.70%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D17" "247.36"
Set-TextValue "E17" "  +8.57%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D18" "63.63"
Set-TextValue "E18" "  +4.58%  "
Set-TextValue "D19" "7.61"
Set-TextValue "E19" "  +3.58%  "
Set-TextValue "E20" "  +3.02%  "
Set-TextValue "E21" "  -0.47%  "
Set-TextValue "D22" "4.06"
Set-TextValue "E22" "  +3.97%  "
Set-TextValue "D23" "9.28"
Set-TextValue "E23" "  +4.10%  "
Set-TextValue "D24" "2.12"
Set-TextValue "E24" "  +4.57%  "
Set-TextValue "D25" "155.88"
Set-TextValue "E25" "  +3.01%  "
Set-TextValue "D26" "15.42"
Set-TextValue "E26" "  +4.71%  "
Set-TextValue "E27" "  +5.92%  "
Set-TextValue "D28" "6.41"
Set-TextValue "E28" "  +2.67%  "
Set-TextValue "D29" "0.995"
Set-TextValue "E29" "  -0.50%  "
Set-TextValue "E30" "  +1.35%  "
Set-TextValue "E31" "  +0.39%  "
Set-TextValue "D32" "3.25"
Set-TextValue "E32" "  +2.62%  "
Set-TextValue "D33" "1.438.44"
Set-TextValue "E33" "  +4.01%  "
Set-TextValue "D34" "3.12"
Set-TextValue "E34" "  +3.74%  "
Set-TextValue "E35" "  -2.25%  "
Set-TextValue "E36" "  +10.80%  "
Set-TextValue "D37" "1.52"
Set-TextValue "E37" "  +2.76%  "
Set-TextValue "D38" "2.30"
Set-TextValue "E38" "  -1.22%  "
Set-TextValue "E39" "  +2.54%  "
Set-TextValue "D40" "0.534"
Set-TextValue "E40" "  +4.52%  "
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D41" "55.11"
Set-TextValue "E41" "  +30.78%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D42" "1.96"
Set-TextValue "E42" "  +1.70%  "
Set-TextValue "D43" "0.801"
Set-TextValue "E43" "  +3.52%  "
Set-TextValue "D44" "0.995"
Set-TextValue "E44" "  -0.49%  "
Set-TextValue "E45" "  +2.73%  "
Set-TextValue "D46" "66.27"
Set-TextValue "E46" "  +6.96%  "
Set-TextValue "D47" "5.33"
Set-TextValue "E47" "  -0.73%  "
Set-TextValue "D48" "1.742.84"
Set-TextValue "E48" "  +3.63%  "
Set-TextValue "D49" "86.33"
Set-TextValue "E49" "  +0.74%  "
Set-TextValue "E50" "  -3.95%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.0519"
Set-TextValue "E51" "  +1.67%  "
